# Add a new worksheet "May 08" after the last existing sheet (May 03),
# matching the target workbook structure (sheetId=4, rId4).
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "May 08"

# --- Header / summary block ---
$ws.Range("A1").Value = "Last Updated"
$ws.Range("B1").Value = "May 08 2022 22:23"

$ws.Range("A3").Value = "Balance"
$ws.Range("B3").Value = -1712.8

$ws.Range("A4").Value = "Total Consumed"
$ws.Range("B4").Value = 1187.2

$ws.Range("A5").Value = "Total Burned"
$ws.Range("B5").Value = 2900

# --- Food table header ---
$ws.Range("A7").Value = "Food"
$ws.Range("B7").Value = "Amount"
$ws.Range("C7").Value = "Unit"
$ws.Range("D7").Value = "Calories"
$ws.Range("E7").Value = "Protein"
$ws.Range("F7").Value = "Carbs"
$ws.Range("G7").Value = "Fats"

# --- Food table rows (numeric-looking values are stored as text, like the
#     source workbook's other day sheets: "1.00", "65.00", etc. The
#     Amount/Calories/Protein/Carbs/Fats columns are formatted as Text so
#     the trailing zeros survive as literal characters.) ---

$ws.Range("B8:G8").NumberFormat = "@"
$ws.Range("A8").Value = "לחם קל (פרוסה)"
$ws.Range("B8").Value = "1.00"
$ws.Range("C8").Value = "יחידות"
$ws.Range("D8").Value = "65.00"
$ws.Range("E8").Value = "4.20"
$ws.Range("F8").Value = "2.60"
$ws.Range("G8").Value = "3.40"

$ws.Range("B9:G9").NumberFormat = "@"
$ws.Range("A9").Value = "טונה במים"
$ws.Range("B9").Value = "1.00"
$ws.Range("C9").Value = "יחידה"
$ws.Range("D9").Value = "120.00"
$ws.Range("E9").Value = "26.00"
$ws.Range("F9").Value = "0.00"
$ws.Range("G9").Value = "1.00"

$ws.Range("B10:G10").NumberFormat = "@"
$ws.Range("A10").Value = "ביצה קשה"
$ws.Range("B10").Value = "2.00"
$ws.Range("C10").Value = "יחידה"
$ws.Range("D10").Value = "180.00"
$ws.Range("E10").Value = "14.60"
$ws.Range("F10").Value = "1.20"
$ws.Range("G10").Value = "12.00"

$ws.Range("B11:G11").NumberFormat = "@"
$ws.Range("A11").Value = "אבוקדו בינוני"
$ws.Range("B11").Value = "0.50"
$ws.Range("C11").Value = "יחידות"
$ws.Range("D11").Value = "140.00"
$ws.Range("E11").Value = "1.75"
$ws.Range("F11").Value = "7.50"
$ws.Range("G11").Value = "13.00"

$ws.Range("B12:G12").NumberFormat = "@"
$ws.Range("A12").Value = "קינואה"
$ws.Range("B12").Value = "200.00"
$ws.Range("C12").Value = "גרם"
$ws.Range("D12").Value = "240.00"
$ws.Range("E12").Value = "8.80"
$ws.Range("F12").Value = "42.60"
$ws.Range("G12").Value = "3.84"

$ws.Range("B13:G13").NumberFormat = "@"
$ws.Range("A13").Value = "חזה עוף 100 גרם"
$ws.Range("B13").Value = "200.00"
$ws.Range("C13").Value = "גרם"
$ws.Range("D13").Value = "228.00"
$ws.Range("E13").Value = "42.00"
$ws.Range("F13").Value = "0.00"
$ws.Range("G13").Value = "5.20"

$ws.Range("B14:G14").NumberFormat = "@"
$ws.Range("A14").Value = "פטל שחור מוקפא"
$ws.Range("B14").Value = "30.00"
$ws.Range("C14").Value = "גרם"
$ws.Range("D14").Value = "17.10"
$ws.Range("E14").Value = "0.51"
$ws.Range("F14").Value = "2.58"
$ws.Range("G14").Value = "0.15"

$ws.Range("B15:G15").NumberFormat = "@"
$ws.Range("A15").Value = "בננה"
$ws.Range("B15").Value = "1.00"
$ws.Range("C15").Value = "יחידות"
$ws.Range("D15").Value = "89.00"
$ws.Range("E15").Value = "1.10"
$ws.Range("F15").Value = "22.80"
$ws.Range("G15").Value = "0.30"

$ws.Range("B16:G16").NumberFormat = "@"
$ws.Range("A16").Value = "סלט ירקות"
$ws.Range("B16").Value = "1.00"
$ws.Range("C16").Value = "יחידה"
$ws.Range("D16").Value = "95.10"
$ws.Range("E16").Value = "4.00"
$ws.Range("F16").Value = "16.80"
$ws.Range("G16").Value = "0.72"

$ws.Range("B17:G17").NumberFormat = "@"
$ws.Range("A17").Value = "מלפפון"
$ws.Range("B17").Value = "1.00"
$ws.Range("C17").Value = "יחידות"
$ws.Range("D17").Value = "13.00"
$ws.Range("E17").Value = "0.70"
$ws.Range("F17").Value = "2.80"
$ws.Range("G17").Value = "0.10"

# Keep the originally active sheet ("May 03") selected/active, since the
# workbook-level active-tab context was not part of the edit.
$wb.Worksheets.Item("May 03").Activate()
